$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = @(
    46073,
    46073.01041666666,
    46073.02083333334,
    46073.03125,
    46073.04166666666,
    46073.05208333334,
    46073.0625,
    46073.07291666666,
    46073.08333333334,
    46073.09375,
    46073.10416666666,
    46073.11458333334,
    46073.125,
    46073.13541666666,
    46073.14583333334,
    46073.15625,
    46073.16666666666,
    46073.17708333334,
    46073.1875,
    46073.19791666666,
    46073.20833333334,
    46073.21875,
    46073.22916666666,
    46073.23958333334,
    46073.25,
    46073.26041666666,
    46073.27083333334,
    46073.28125,
    46073.29166666666,
    46073.30208333334,
    46073.3125,
    46073.32291666666,
    46073.33333333334,
    46073.34375,
    46073.35416666666,
    46073.36458333334,
    46073.375,
    46073.38541666666,
    46073.39583333334,
    46073.40625,
    46073.41666666666,
    46073.42708333334,
    46073.4375,
    46073.44791666666,
    46073.45833333334,
    46073.46875,
    46073.47916666666,
    46073.48958333334,
    46073.5,
    46073.51041666666,
    46073.52083333334,
    46073.53125,
    46073.54166666666,
    46073.55208333334,
    46073.5625,
    46073.57291666666,
    46073.58333333334,
    46073.59375,
    46073.60416666666,
    46073.61458333334,
    46073.625,
    46073.63541666666,
    46073.64583333334,
    46073.65625,
    46073.66666666666,
    46073.67708333334,
    46073.6875,
    46073.69791666666,
    46073.70833333334,
    46073.71875,
    46073.72916666666,
    46073.73958333334,
    46073.75,
    46073.76041666666,
    46073.77083333334,
    46073.78125,
    46073.79166666666,
    46073.80208333334,
    46073.8125,
    46073.82291666666,
    46073.83333333334,
    46073.84375,
    46073.85416666666,
    46073.86458333334,
    46073.875,
    46073.88541666666,
    46073.89583333334,
    46073.90625,
    46073.91666666666,
    46073.92708333334,
    46073.9375,
    46073.94791666666,
    46073.95833333334,
    46073.96875,
    46073.97916666666,
    46073.98958333334,
    46074,
    46074.01041666666,
    46074.02083333334,
    46074.03125,
    46074.04166666666,
    46074.05208333334,
    46074.0625,
    46074.07291666666,
    46074.08333333334,
    46074.09375,
    46074.10416666666,
    46074.11458333334,
    46074.125,
    46074.13541666666,
    46074.14583333334,
    46074.15625,
    46074.16666666666,
    46074.17708333334,
    46074.1875,
    46074.19791666666,
    46074.20833333334,
    46074.21875,
    46074.22916666666,
    46074.23958333334,
    46074.25,
    46074.26041666666,
    46074.27083333334,
    46074.28125,
    46074.29166666666,
    46074.30208333334,
    46074.3125,
    46074.32291666666,
    46074.33333333334,
    46074.34375,
    46074.35416666666,
    46074.36458333334,
    46074.375,
    46074.38541666666,
    46074.39583333334,
    46074.40625,
    46074.41666666666,
    46074.42708333334,
    46074.4375,
    46074.44791666666,
    46074.45833333334,
    46074.46875,
    46074.47916666666,
    46074.48958333334,
    46074.5,
    46074.51041666666,
    46074.52083333334,
    46074.53125,
    46074.54166666666,
    46074.55208333334,
    46074.5625,
    46074.57291666666,
    46074.58333333334,
    46074.59375,
    46074.60416666666,
    46074.61458333334,
    46074.625,
    46074.63541666666,
    46074.64583333334,
    46074.65625,
    46074.66666666666,
    46074.67708333334,
    46074.6875,
    46074.69791666666,
    46074.70833333334,
    46074.71875,
    46074.72916666666,
    46074.73958333334,
    46074.75,
    46074.76041666666,
    46074.77083333334,
    46074.78125,
    46074.79166666666,
    46074.80208333334,
    46074.8125,
    46074.82291666666,
    46074.83333333334,
    46074.84375,
    46074.85416666666,
    46074.86458333334,
    46074.875,
    46074.88541666666,
    46074.89583333334,
    46074.90625,
    46074.91666666666,
    46074.92708333334,
    46074.9375,
    46074.94791666666,
    46074.95833333334,
    46074.96875,
    46074.97916666666,
    46074.98958333334
)

$colB = @(
    1640.767,
    1633,
    1615.019,
    1594.487,
    1588.394,
    1561.678,
    1534.785,
    1498.41,
    1438.709,
    1404.629,
    1374.104,
    1345.101,
    1279.825,
    1246.97,
    1241.242,
    1225.905,
    1148.853,
    1122.907,
    1099.461,
    1075.97,
    970.7569999999999,
    944.0940000000001,
    933.802,
    907.455,
    855.3440000000001,
    819.669,
    738.621,
    716.359,
    670.648,
    642.188,
    618.982,
    601.327,
    578.015,
    566.846,
    551.5940000000001,
    531.289,
    572.735,
    562.938,
    554.414,
    545.407,
    561.799,
    569.877,
    579.498,
    588.889,
    605.377,
    634.264,
    664.782,
    695.9930000000001,
    765.553,
    817.776,
    871.354,
    925.99,
    1044.103,
    1135.436,
    1204.974,
    1319.811,
    1440.112,
    1527.69,
    1568.443,
    1647.216,
    1735.951,
    1794.899,
    1904.164,
    1963.702,
    2054.848,
    2095.504,
    2136.05,
    2176.247,
    2225.217,
    2246.77,
    2268.106,
    2288.36,
    2319.606,
    2326.041,
    2331.66,
    2336.534,
    2336.843,
    2335.701,
    2335.241,
    2334.972,
    2341.798,
    2332.354,
    2322.34,
    2250.323,
    2306.569,
    2243.223,
    2239.162,
    2237.393,
    2226.751,
    2221.989,
    2217.946,
    2213.739,
    2300.605,
    2297.812,
    2294.392,
    2291.212,
    2302.824,
    2290.375,
    2283.566,
    2277.197,
    2273.998,
    2258.526,
    2256.114,
    2240.242,
    2232.002,
    2228.131,
    2224.797,
    2220.095,
    2222.537,
    2214.861,
    2207.065,
    2206.985,
    2187.612,
    2180.397,
    2158.825,
    2152.32,
    2122.388,
    2111.598,
    2108.99,
    2098.661,
    2060.225,
    2060.313,
    2051.123,
    2042.162,
    2015.771,
    2005.438,
    1995.037,
    1984.012,
    1953.954,
    1938.69,
    1929.549,
    1916.242,
    1898.887,
    1888.365,
    1878.002,
    1867.168,
    1853.442,
    1842.567,
    1831.311,
    1810.691,
    1795.835,
    1778.907,
    1760.429,
    1733.07,
    1699.638,
    1670.644,
    1650.522,
    1630.916,
    1599.797,
    1583.454,
    1559.865,
    1545.736,
    1508.666,
    1489.409,
    1478.327,
    1458,
    1427.209,
    1417.418,
    1399.506,
    1381.071,
    1381.851,
    1363.691,
    1383.396,
    1364.168,
    1329.668,
    1309.73,
    1290.437,
    1271.469,
    1237.685,
    1212.919,
    1186.775,
    1162.391,
    1121.232,
    1096.588,
    1072.018,
    1047.866,
    1015.289,
    994.963,
    973.862,
    933.394,
    903.005,
    884.022,
    866.254,
    848.0700000000001,
    831.002,
    825.068,
    817.423,
    810.174,
    0,
    0,
    0,
    0
)

$colC = @(
    1643,
    1561,
    1487,
    1402,
    1377,
    1365,
    1328,
    1314,
    1380,
    1397,
    1396,
    1341,
    1307,
    1305,
    1272,
    1175,
    1060,
    998,
    942,
    815,
    678,
    568,
    503,
    456,
    425,
    384,
    364,
    349,
    326,
    297,
    265,
    260,
    270,
    288,
    308,
    321,
    333,
    323,
    317,
    318,
    352,
    386,
    433,
    472,
    520,
    591,
    637,
    680,
    712,
    748,
    808,
    936,
    1033,
    1090,
    1168,
    1306,
    1509,
    1627,
    1733,
    1821,
    1869,
    1870,
    2038,
    2099,
    2168,
    2189,
    2231,
    2270,
    2355,
    2387,
    2423,
    2458,
    2519,
    2532,
    2540,
    2519,
    2518,
    2504,
    0,
    2477,
    2479,
    2465,
    2446,
    2362,
    2400,
    2334,
    2293,
    2295,
    2386,
    2411,
    2355,
    2378,
    2445,
    2458,
    2455,
    2460,
    2475,
    2473,
    2443,
    2415,
    2385,
    2373,
    2382,
    2392,
    0,
    2386,
    2381,
    2370,
    2350,
    2339,
    2296,
    2259,
    2238,
    2185,
    2094,
    1956,
    1877,
    1851,
    1847,
    1812,
    1676,
    1692,
    1651,
    1611,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0
)

$colE = @(
    "20.02.20261",
    "20.02.20262",
    "20.02.20263",
    "20.02.20264",
    "20.02.20265",
    "20.02.20266",
    "20.02.20267",
    "20.02.20268",
    "20.02.20269",
    "20.02.202610",
    "20.02.202611",
    "20.02.202612",
    "20.02.202613",
    "20.02.202614",
    "20.02.202615",
    "20.02.202616",
    "20.02.202617",
    "20.02.202618",
    "20.02.202619",
    "20.02.202620",
    "20.02.202621",
    "20.02.202622",
    "20.02.202623",
    "20.02.202624",
    "20.02.202625",
    "20.02.202626",
    "20.02.202627",
    "20.02.202628",
    "20.02.202629",
    "20.02.202630",
    "20.02.202631",
    "20.02.202632",
    "20.02.202633",
    "20.02.202634",
    "20.02.202635",
    "20.02.202636",
    "20.02.202637",
    "20.02.202638",
    "20.02.202639",
    "20.02.202640",
    "20.02.202641",
    "20.02.202642",
    "20.02.202643",
    "20.02.202644",
    "20.02.202645",
    "20.02.202646",
    "20.02.202647",
    "20.02.202648",
    "20.02.202649",
    "20.02.202650",
    "20.02.202651",
    "20.02.202652",
    "20.02.202653",
    "20.02.202654",
    "20.02.202655",
    "20.02.202656",
    "20.02.202657",
    "20.02.202658",
    "20.02.202659",
    "20.02.202660",
    "20.02.202661",
    "20.02.202662",
    "20.02.202663",
    "20.02.202664",
    "20.02.202665",
    "20.02.202666",
    "20.02.202667",
    "20.02.202668",
    "20.02.202669",
    "20.02.202670",
    "20.02.202671",
    "20.02.202672",
    "20.02.202673",
    "20.02.202674",
    "20.02.202675",
    "20.02.202676",
    "20.02.202677",
    "20.02.202678",
    "20.02.202679",
    "20.02.202680",
    "20.02.202681",
    "20.02.202682",
    "20.02.202683",
    "20.02.202684",
    "20.02.202685",
    "20.02.202686",
    "20.02.202687",
    "20.02.202688",
    "20.02.202689",
    "20.02.202690",
    "20.02.202691",
    "20.02.202692",
    "20.02.202693",
    "20.02.202694",
    "20.02.202695",
    "20.02.202696",
    "21.02.20261",
    "21.02.20262",
    "21.02.20263",
    "21.02.20264",
    "21.02.20265",
    "21.02.20266",
    "21.02.20267",
    "21.02.20268",
    "21.02.20269",
    "21.02.202610",
    "21.02.202611",
    "21.02.202612",
    "21.02.202613",
    "21.02.202614",
    "21.02.202615",
    "21.02.202616",
    "21.02.202617",
    "21.02.202618",
    "21.02.202619",
    "21.02.202620",
    "21.02.202621",
    "21.02.202622",
    "21.02.202623",
    "21.02.202624",
    "21.02.202625",
    "21.02.202626",
    "21.02.202627",
    "21.02.202628",
    "21.02.202629",
    "21.02.202630",
    "21.02.202631",
    "21.02.202632",
    "21.02.202633",
    "21.02.202634",
    "21.02.202635",
    "21.02.202636",
    "21.02.202637",
    "21.02.202638",
    "21.02.202639",
    "21.02.202640",
    "21.02.202641",
    "21.02.202642",
    "21.02.202643",
    "21.02.202644",
    "21.02.202645",
    "21.02.202646",
    "21.02.202647",
    "21.02.202648",
    "21.02.202649",
    "21.02.202650",
    "21.02.202651",
    "21.02.202652",
    "21.02.202653",
    "21.02.202654",
    "21.02.202655",
    "21.02.202656",
    "21.02.202657",
    "21.02.202658",
    "21.02.202659",
    "21.02.202660",
    "21.02.202661",
    "21.02.202662",
    "21.02.202663",
    "21.02.202664",
    "21.02.202665",
    "21.02.202666",
    "21.02.202667",
    "21.02.202668",
    "21.02.202669",
    "21.02.202670",
    "21.02.202671",
    "21.02.202672",
    "21.02.202673",
    "21.02.202674",
    "21.02.202675",
    "21.02.202676",
    "21.02.202677",
    "21.02.202678",
    "21.02.202679",
    "21.02.202680",
    "21.02.202681",
    "21.02.202682",
    "21.02.202683",
    "21.02.202684",
    "21.02.202685",
    "21.02.202686",
    "21.02.202687",
    "21.02.202688",
    "21.02.202689",
    "21.02.202690",
    "21.02.202691",
    "21.02.202692",
    "21.02.202693",
    "21.02.202694",
    "21.02.202695",
    "21.02.202696"
)

$startRow = 2
for ($i = 0; $i -lt $colA.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value2 = $colA[$i]
    $ws.Cells.Item($r, 2).Value2 = $colB[$i]
    $ws.Cells.Item($r, 3).Value2 = $colC[$i]
    $ws.Cells.Item($r, 5).Value2 = $colE[$i]
}
